$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 13: E13 value change
$ws.Range("E13").Value = 332114255.60000002

# Row 14: E14 value change
$ws.Range("E14").Value = 537958.22

# Row 15: E15 value change
$ws.Range("E15").Value = 30000000

# Row 16: E16 value change
$ws.Range("E16").Value = -60473972.810000002

# Row 18: E18 becomes a formula =SUM(E12:E17)
$ws.Range("E18").Formula = "=SUM(E12:E17)"

# Row 19: E19 value change
$ws.Range("E19").Value = -384700000

# Row 21: E21 becomes a formula =SUM(E18:E20)
$ws.Range("E21").Formula = "=SUM(E18:E20)"

# Row 22: E22 value change
$ws.Range("E22").Value = -20015625

$wb.Save()
